$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a cell's value as literal text while preserving its original
# style / number format. This stops Excel's automatic type inference from
# turning numeric-looking strings (e.g. "243.16") into real numbers.
function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '30.286.13'
Set-TextValue $ws.Range("E2") '  -1.17%  '
Set-TextValue $ws.Range("D3") '1.865.33'
Set-TextValue $ws.Range("E3") '  -0.93%  '
Set-TextValue $ws.Range("E4") '  +0.00%  '
Set-TextValue $ws.Range("D5") '243.16'
Set-TextValue $ws.Range("D7") '0.4723'
Set-TextValue $ws.Range("E7") '  -0.60%  '
Set-TextValue $ws.Range("D8") '0.2871'
Set-TextValue $ws.Range("E8") '  -2.37%  '
Set-TextValue $ws.Range("D9") '0.06469'
Set-TextValue $ws.Range("E9") '  -1.07%  '
Set-TextValue $ws.Range("D10") '21.52'
Set-TextValue $ws.Range("E10") '  -2.14%  '
Set-TextValue $ws.Range("D11") '0.07788'
Set-TextValue $ws.Range("E11") '  +0.53%  '
Set-TextValue $ws.Range("D12") '96.51'
Set-TextValue $ws.Range("E12") '  -0.42%  '
Set-TextValue $ws.Range("D13") '1.866.22'
Set-TextValue $ws.Range("E13") '  -0.83%  '
Set-TextValue $ws.Range("D14") '0.7186'
Set-TextValue $ws.Range("E14") '  -2.57%  '
Set-TextValue $ws.Range("D15") '5.124'
Set-TextValue $ws.Range("E15") '  -2.29%  '
Set-TextValue $ws.Range("D16") '279.01'
Set-TextValue $ws.Range("E16") '  +1.43%  '
Set-TextValue $ws.Range("D17") '30.276.45'
Set-TextValue $ws.Range("E17") '  -1.12%  '
Set-TextValue $ws.Range("D18") '12.97'
Set-TextValue $ws.Range("E18") '  -1.64%  '
Set-TextValue $ws.Range("E19") '  -0.03%  '
Set-TextValue $ws.Range("D20") '0.000007457'
Set-TextValue $ws.Range("E20") '  -1.09%  '
Set-TextValue $ws.Range("D21") '2.108.97'
Set-TextValue $ws.Range("E21") '  -0.86%  '
Set-TextValue $ws.Range("E22") '  -0.02%  '
Set-TextValue $ws.Range("D23") '5.230'
Set-TextValue $ws.Range("E23") '  -2.17%  '
Set-TextValue $ws.Range("D24") '6.235'
Set-TextValue $ws.Range("E24") '  -0.02%  '
Set-TextValue $ws.Range("D25") '161.89'
Set-TextValue $ws.Range("E25") '  -1.27%  '
Set-TextValue $ws.Range("D26") '8.956'
Set-TextValue $ws.Range("E26") '  -2.97%  '
Set-TextValue $ws.Range("D27") '18.64'
Set-TextValue $ws.Range("E27") '  -1.16%  '
Set-TextValue $ws.Range("D28") '1.873'
Set-TextValue $ws.Range("E28") '  -1.98%  '
Set-TextValue $ws.Range("D29") '0.09619'
Set-TextValue $ws.Range("E29") '  -1.11%  '
Set-TextValue $ws.Range("D30") '1.312'
Set-TextValue $ws.Range("E30") '  -2.43%  '
Set-TextValue $ws.Range("D31") '1.477'
Set-TextValue $ws.Range("E31") '  -1.90%  '
Set-TextValue $ws.Range("D32") '4.205'
Set-TextValue $ws.Range("E32") '  -2.04%  '
Set-TextValue $ws.Range("D33") '4.103'
Set-TextValue $ws.Range("E33") '  -1.29%  '
Set-TextValue $ws.Range("D34") '0.04765'
Set-TextValue $ws.Range("E34") '  -2.09%  '
Set-TextValue $ws.Range("D35") '1.115'
Set-TextValue $ws.Range("E35") '  -0.96%  '
Set-TextValue $ws.Range("D36") '0.6820'
Set-TextValue $ws.Range("E36") '  -2.57%  '
Set-TextValue $ws.Range("E37") '  -0.36%  '
Set-TextValue $ws.Range("D38") '0.01887'
Set-TextValue $ws.Range("E38") '  -1.15%  '
Set-TextValue $ws.Range("D39") '2.838'
Set-TextValue $ws.Range("E39") '  +1.70%  '
Set-TextValue $ws.Range("D40") '75.18'
Set-TextValue $ws.Range("E40") '  +0.28%  '
Set-TextValue $ws.Range("D41") '6.184'
Set-TextValue $ws.Range("E41") '  -1.96%  '
Set-TextValue $ws.Range("D42") '1.924'
Set-TextValue $ws.Range("E42") '  -5.30%  '
Set-TextValue $ws.Range("D43") '0.4191'
Set-TextValue $ws.Range("E43") '  -1.45%  '
Set-TextValue $ws.Range("D44") '0.9993'
Set-TextValue $ws.Range("E44") '  -0.09%  '
Set-TextValue $ws.Range("D45") '0.8245'
Set-TextValue $ws.Range("E45") '  -2.03%  '
Set-TextValue $ws.Range("D46") '100.33'
Set-TextValue $ws.Range("E46") '  -2.22%  '
Set-TextValue $ws.Range("D47") '9.563'
Set-TextValue $ws.Range("E47") '  +1.32%  '
Set-TextValue $ws.Range("D48") '6.938'
Set-TextValue $ws.Range("E48") '  -1.57%  '
Set-TextValue $ws.Range("D49") '34.89'
Set-TextValue $ws.Range("E49") '  -2.06%  '
Set-TextValue $ws.Range("D50") '0.05768'
Set-TextValue $ws.Range("E50") '  -0.06%  '
Set-TextValue $ws.Range("D51") '881.91'
Set-TextValue $ws.Range("E51") '  -3.89%  '
